$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 18184118
$ws.Range("I116").Value = 33335284
$ws.Range("J116").Value = 2719
$ws.Range("K116").Value = 33335284
$ws.Range("L116").Value = 2719
$ws.Range("M116").Value = -33331842
$ws.Range("N116").Value = -9603
$ws.Range("H132").Value = 1836.7091
$ws.Range("I132").Value = 1677.0605
$ws.Range("J132").Value = 2076.182
$ws.Range("K132").Value = 5031.181500000001
$ws.Range("L132").Value = 6228.545999999999
$ws.Range("M132").Value = -2501.181500000001
$ws.Range("N132").Value = -11288.546
$ws.Range("H137").Value = 1189.2565
$ws.Range("I137").Value = 1163.6571
$ws.Range("J137").Value = 1413.25
$ws.Range("K137").Value = 3490.9713
$ws.Range("L137").Value = 4239.75
$ws.Range("M137").Value = -940.9712999999997
$ws.Range("N137").Value = -9339.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2316.4092
$ws.Range("I61").Value = 2115.353
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2115.353
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1903.353
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1115.8235
$ws.Range("I74").Value = 830.75
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 830.75
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = 43.25
$ws.Range("N74").Value = -3548
$ws.Range("H77").Value = 1115.8235
$ws.Range("I77").Value = 830.75
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 4153.75
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = 214.25
$ws.Range("N77").Value = -17736
$ws.Range("H123").Value = 24174.334
$ws.Range("J123").Value = 24174.334
$ws.Range("L123").Value = 24174.334
$ws.Range("N123").Value = -33974.334
$ws.Range("H132").Value = 3715.0908
$ws.Range("I132").Value = 4099.737
$ws.Range("J132").Value = 2855.2942
$ws.Range("K132").Value = 12299.211
$ws.Range("L132").Value = 8565.882599999999
$ws.Range("M132").Value = -9769.210999999999
$ws.Range("N132").Value = -13625.8826
$ws.Range("H136").Value = 2316.4092
$ws.Range("I136").Value = 2115.353
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6346.059
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3796.059
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2283.0476
$ws.Range("I134").Value = 1951
$ws.Range("J134").Value = 3023.7693
$ws.Range("K134").Value = 5853
$ws.Range("L134").Value = 9071.3079
$ws.Range("M134").Value = -3318
$ws.Range("N134").Value = -14141.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1838.3549
$ws.Range("I31").Value = 1481.1111
$ws.Range("J31").Value = 4249.75
$ws.Range("K31").Value = 1481.1111
$ws.Range("L31").Value = 4249.75
$ws.Range("M31").Value = -1186.1111
$ws.Range("N31").Value = -4839.75
$ws.Range("H34").Value = 1838.3549
$ws.Range("I34").Value = 1481.1111
$ws.Range("J34").Value = 4249.75
$ws.Range("K34").Value = 1481.1111
$ws.Range("L34").Value = 4249.75
$ws.Range("M34").Value = -1279.1111
$ws.Range("N34").Value = -4653.75
$ws.Range("H58").Value = 773328.9
$ws.Range("I58").Value = 1278727.4
$ws.Range("J58").Value = 1931.0526
$ws.Range("K58").Value = 1278727.4
$ws.Range("L58").Value = 1931.0526
$ws.Range("M58").Value = -1278524.4
$ws.Range("N58").Value = -2337.0526
$ws.Range("H132").Value = 323350.66
$ws.Range("I132").Value = 410761.06
$ws.Range("J132").Value = 2845.7778
$ws.Range("K132").Value = 1232283.18
$ws.Range("L132").Value = 8537.3334
$ws.Range("M132").Value = -1229753.18
$ws.Range("N132").Value = -13597.3334
$ws.Range("H134").Value = 1671.6111
$ws.Range("I134").Value = 1197.1052
$ws.Range("J134").Value = 2798.5625
$ws.Range("K134").Value = 3591.3156
$ws.Range("L134").Value = 8395.6875
$ws.Range("M134").Value = -1056.3156
$ws.Range("N134").Value = -13465.6875
$ws.Range("H136").Value = 773328.9
$ws.Range("I136").Value = 1278727.4
$ws.Range("J136").Value = 1931.0526
$ws.Range("K136").Value = 3836182.2
$ws.Range("L136").Value = 5793.1578
$ws.Range("M136").Value = -3833632.2
$ws.Range("N136").Value = -10893.1578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 131.9375
$ws.Range("I38").Value = 150.9
$ws.Range("J38").Value = 100.333336
$ws.Range("K38").Value = 452.7
$ws.Range("L38").Value = 301.000008
$ws.Range("M38").Value = -105.7
$ws.Range("N38").Value = -995.000008
$ws.Range("H131").Value = 1160.5294
$ws.Range("J131").Value = 1292
$ws.Range("L131").Value = 3876
$ws.Range("N131").Value = -13956

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3842.182
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4158
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12474
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -17374
$ws.Range("H126").Value = 3951.8
$ws.Range("I126").Value = 3889.0667
$ws.Range("J126").Value = 4140
$ws.Range("K126").Value = 11667.2001
$ws.Range("L126").Value = 12420
$ws.Range("M126").Value = -9197.2001
$ws.Range("N126").Value = -17360
$ws.Range("H132").Value = 2612.0715
$ws.Range("I132").Value = 1941.4445
$ws.Range("J132").Value = 3819.2
$ws.Range("K132").Value = 5824.333500000001
$ws.Range("L132").Value = 11457.6
$ws.Range("M132").Value = -3294.333500000001
$ws.Range("N132").Value = -16517.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6793.75
$ws.Range("I40").Value = 8125
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 8125
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -7989
$ws.Range("N40").Value = -3072
$ws.Range("H46").Value = 876
$ws.Range("I46").Value = 657.8946999999999
$ws.Range("J46").Value = 1566.6666
$ws.Range("K46").Value = 657.8946999999999
$ws.Range("L46").Value = 1566.6666
$ws.Range("M46").Value = -469.8946999999999
$ws.Range("N46").Value = -1942.6666
$ws.Range("H132").Value = 5314.9
$ws.Range("I132").Value = 6023.3687
$ws.Range("J132").Value = 4091.182
$ws.Range("K132").Value = 18070.1061
$ws.Range("L132").Value = 12273.546
$ws.Range("M132").Value = -15540.1061
$ws.Range("N132").Value = -17333.546
$ws.Range("H136").Value = 29870800
$ws.Range("I136").Value = 38794810
$ws.Range("K136").Value = 116384430
$ws.Range("M136").Value = -116381880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 14933.333
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 14933.333
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 14933.333
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -16077.333
$ws.Range("H113").Value = 328.92856
$ws.Range("I113").Value = 445.33334
$ws.Range("J113").Value = 119.4
$ws.Range("K113").Value = 1336.00002
$ws.Range("L113").Value = 358.2
$ws.Range("M113").Value = 833.9999800000001
$ws.Range("N113").Value = -4698.2
$ws.Range("H126").Value = 3782.7896
$ws.Range("I126").Value = 3955.5715
$ws.Range("J126").Value = 3299
$ws.Range("K126").Value = 11866.7145
$ws.Range("L126").Value = 9897
$ws.Range("M126").Value = -9396.7145
$ws.Range("N126").Value = -14837
$ws.Range("H132").Value = 2107.5247
$ws.Range("I132").Value = 1279.3889
$ws.Range("J132").Value = 3300.04
$ws.Range("K132").Value = 3838.1667
$ws.Range("L132").Value = 9900.119999999999
$ws.Range("M132").Value = -1308.1667
$ws.Range("N132").Value = -14960.12
$ws.Range("H136").Value = 1612.6285
$ws.Range("I136").Value = 1544.1786
$ws.Range("J136").Value = 1886.4286
$ws.Range("K136").Value = 4632.5358
$ws.Range("L136").Value = 5659.2858
$ws.Range("M136").Value = -2082.5358
$ws.Range("N136").Value = -10759.2858
